$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Color, Size, Style, Supplier values in row 2
$ws.Range("I2").Value = "L GREY"

# Size and Style columns hold numeric-looking text; force text storage
# then restore the default style so no stray formatting is left behind.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "36"
$ws.Range("J2").Style = "Normal"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "4444"
$ws.Range("L2").Style = "Normal"

$ws.Range("N2").Value = "Supplier4"

# Update Description (total wsp) numeric value
$ws.Range("G2").Value = 41
